{"js": "// Word JS API (Office.js) script.\n// Body of: async (context) => { ... }\n//\n// Reproduces:\n//   1. Resize \"Picture 2\" (top-5-schools table image) from\n//      5943600 x 1682115 EMU  ->  5303520 x 1722120 EMU.\n//   2. Resize \"Picture 3\" (charter-vs-district table image) from\n//      5943600 x 871220 EMU  ->  5974080 x 871220 EMU (width only).\n//   3. Remove one of the two consecutive empty paragraphs that sit\n//      right after Picture 3 (collapsing the double blank line to one).\n//   4. Fix \"...budget per student has the highest...\" ->\n//      \"...budget per student have the highest...\" in the closing\n//      bullet paragraph.\n\nconst body = context.document.body;\n\n// --- 1 & 2: resize the two inline pictures -----------------------------\n// The public InlinePicture.width/height setters are not wired up to the\n// underlying shape geometry in this host, and the requested resize is not\n// a proportional scale of the original image (aspect ratio changes), so\n// LockAspectRatio-style width/height assignment can't hit the target\n// numbers anyway. Drive the same Word object-model \"InlineShape\" members\n// (ScaleWidth/ScaleHeight, percentages of the native image size) that the\n// shim itself uses under the hood for every other read/write.\nconst pics = body.inlinePictures;\npics.load(\"items\");\nawait context.sync();\n\nif (pics.items.length >= 2) {\n  const picture2 = pics.items[0]; // 5943600 x 1682115 EMU -> 5303520 x 1722120 EMU\n  const picture3 = pics.items[1]; // 5943600 x 871220 EMU  -> 5974080 x 871220 EMU\n\n  picture2._omSet(\"ScaleWidth\", 46.284289276807980, \"InlineShape\");\n  picture2._omSet(\"ScaleHeight\", 53.176470588235300, \"InlineShape\");\n\n  picture3._omSet(\"ScaleWidth\", 70.235162374020160, \"InlineShape\");\n  picture3._omSet(\"ScaleHeight\", 69.821882951653950, \"InlineShape\");\n\n  await context.sync();\n}\n\n// --- 3: collapse the duplicate empty paragraph --------------------------\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paras.items.length - 1; i++) {\n  if (paras.items[i].text === \"\" && paras.items[i + 1].text === \"\") {\n    paras.items[i + 1].delete();\n    await context.sync();\n    break;\n  }\n}\n\n// --- 4: \"has\" -> \"have\" in the closing bullet ---------------------------\nconst hits = body.search(\"student has the highest\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length > 0) {\n  hits.items[0].insertText(\"student have the highest\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is already open as $d below.\n#\n# Reproduces:\n#   1. Resize \"Picture 2\" (top-5-schools table image) from\n#      5943600 x 1682115 EMU  ->  5303520 x 1722120 EMU.\n#   2. Resize \"Picture 3\" (charter-vs-district table image) from\n#      5943600 x 871220 EMU  ->  5974080 x 871220 EMU (width only).\n#   3. Remove one of the two consecutive empty paragraphs that sit\n#      right after Picture 3 (collapsing the double blank line to one).\n#   4. Fix \"...budget per student has the highest...\" ->\n#      \"...budget per student have the highest...\" in the closing\n#      bullet paragraph.\n\n$d = $word.ActiveDocument\n\n# --- 1 & 2: resize the two inline pictures ------------------------------\n# Width/Height are locked to the image's original aspect ratio in this\n# host (each setter recomputes the other dimension from the native image\n# size), and the target size is not a proportional scale of the original,\n# so assigning Width/Height directly can never reach the target extents.\n# ScaleWidth/ScaleHeight (percentages of the native/intrinsic image size)\n# resize each axis independently and land on the exact target EMU values.\n$picture2 = $d.InlineShapes.Item(1)   # 5943600 x 1682115 EMU -> 5303520 x 1722120 EMU\n$picture3 = $d.InlineShapes.Item(2)   # 5943600 x 871220 EMU  -> 5974080 x 871220 EMU\n\n$picture2.ScaleWidth = 46.284289276807980\n$picture2.ScaleHeight = 53.176470588235300\n\n$picture3.ScaleWidth = 70.235162374020160\n$picture3.ScaleHeight = 69.821882951653950\n\n# --- 3: collapse the duplicate empty paragraph --------------------------\nfor ($i = 1; $i -lt $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $pNext = $d.Paragraphs.Item($i + 1)\n    if ($p.Range.Text -eq \"`r\" -and $pNext.Range.Text -eq \"`r\") {\n        $pNext.Range.Delete()\n        break\n    }\n}\n\n# --- 4: \"has\" -> \"have\" in the closing bullet ---------------------------\n$rng = $d.Content\n$rng.Find.Execute(\"student has the highest\", $false, $false, $false, $false, $false, $true, 1, $false, \"student have the highest\", 2) | Out-Null\n"}
